# Fixed issue with Analytics not working
# Update the day's agent-handoff record in row 2 with corrected values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 354
$ws.Range("D2").Value = "Dawnetta"
$ws.Range("E2").Value = "ADAM"
$ws.Range("F2").Value = "|| 17:32-EDT | 04/12/2023 ||"
$ws.Range("G2").Value = "ANTHONY"
